$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily rows (date serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila ab.)
# aggiornamento fino al 6 gennaio 2022
$data = @(
    @(44539,0,0,0),
    @(44540,0,0,0),
    @(44541,0,0,0),
    @(44542,0,0,0),
    @(44543,0,0,0),
    @(44544,0,0,0),
    @(44545,0,0,0),
    @(44546,0,0,0),
    @(44547,1,1,62.34413965087282),
    @(44548,0,1,62.34413965087282),
    @(44550,1,2,124.6882793017456),
    @(44551,1,3,187.0324189526185),
    @(44552,0,3,187.0324189526185),
    @(44553,0,3,187.0324189526185),
    @(44554,0,3,187.0324189526185),
    @(44555,0,2,124.6882793017456),
    @(44556,1,3,187.0324189526185),
    @(44557,2,4,249.3765586034913),
    @(44558,0,3,187.0324189526185),
    @(44559,0,3,187.0324189526185),
    @(44560,1,4,249.3765586034913),
    @(44561,0,4,249.3765586034913),
    @(44562,2,6,374.0648379052369),
    @(44563,0,5,311.7206982543641),
    @(44564,1,4,249.3765586034913),
    @(44565,0,4,249.3765586034913),
    @(44566,4,8,498.7531172069826)
)

$startRow = 465

# Write the new block of values
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

$endRow = $startRow + $data.Count - 1

# Column A carries the date style (same as the existing rows above) - copy
# the formatting from the last pre-existing date cell so the new cells
# reuse the workbook's existing style record instead of minting a new one.
$ws.Range("A464").Copy() | Out-Null
$ws.Range("A$startRow`:A$endRow").PasteSpecial(-4122) | Out-Null
